$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
# D1 was "Publish the Application" slot's sibling index; value stays "Time [h]"
# but the shared-string table reflows once "Publish the Application" is dropped.
$ws.Range("D1").Value = "Time [h]"

# --- Fix hours / task text across the existing rows ---
$ws.Range("D3").Value = 3

$ws.Range("C5").Value = "Implement a DataGridView for Displaying the Data pulled from the DB"

$ws.Range("D10").Value = 9

$ws.Range("C11").Value = "Implement onClick generation of a new Window with the detailed info of the clicked member in the table."

$ws.Range("C12").Value = "Build a member info (new payment) update section in the detailed information Window"

$ws.Range("C13").Value = "Implement the addition of a new Member in the Database via the New Member Window"
$ws.Range("D13").Value = 3.5

$ws.Range("C14").Value = "Implement the BackupManager, making a Backup before each change in the DB."
$ws.Range("D14").Value = 4

# --- Add the new Sprint 5 task row (row 15), matching formatting of row 14 ---
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A15").Value = "Sprint"
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "Implement removal of members form the DB"
$ws.Range("D15").Value = 4.5

# --- Drop the old totals row 16 and rebuild it at row 17 to include the new row ---
$ws.Range("A16:D16").EntireRow.Delete()
$ws.Range("D15").Copy()
$ws.Range("D17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D17").Formula = "=SUM(D2:D15)"

# --- View adjustments ---
$ws.Range("D11").Select()
$excel.ActiveWindow.ScrollColumn = 2
